$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "seq_name" column (B) values were re-cased from "revA_N" to "RevA_N"
# for rows 2 through 37 (revA_1 .. revA_36 -> RevA_1 .. RevA_36).
for ($i = 1; $i -le 36; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 2).Value = "RevA_$i"
}

# Reflect the author's final selection in the saved view (B2:B37, active cell B2).
$ws.Range("B2:B37").Select()
